# Debug why the password verifier is not checking credentials in real time
#
# This script mimics a user editing the "students" workbook by hand:
# overwriting some test values on existing rows and then adding a batch
# of dummy student sign-up rows (and matching password rows) to exercise
# the real-time password verifier.

$wb = $excel.ActiveWorkbook

$students = $wb.Worksheets.Item("students")
$pswd     = $wb.Worksheets.Item("student_pswd")

# --- overwrite some existing junk/test values on the "students" sheet ---
$students.Range("C2").Value = "hj"
$students.Range("D2").Value = "hjk"
$students.Range("E2").Value = "hj"
$students.Range("F2").Value = "kh"
$students.Range("G2").Value = "jkhhj"
$students.Range("H2").Value = "hjkh"

$students.Range("C3").Value = "fdfd"
$students.Range("D3").Value = "g"
$students.Range("E3").Value = "hg"
$students.Range("F3").Value = "hg"
$students.Range("G3").Value = "hjg"
$students.Range("H3").Value = "hgj"
$students.Range("J3").Value = 8

# --- add new dummy students (rows 4-9) on "students" ---
$newStudents = @(
    @(3, "Sam",     "fj",    "hj",      "h",       "jkh",   "jk",  "h"),
    @(4, "Agnetha", "hjk",   "h",       "jh",      "jkh",   "jk",  "hj"),
    @(5, "Lee",     "dffd",  "hj",      "hj",      "kh",    "jk",  "hj"),
    @(6, "Garry",   "fgfgf", "h",       "jh",      "jhj",   "kh",  "jk"),
    @(7, "Kumara",  "hjkh",  "jkhjh",   "hjkhjkh", "hjkhh", "hjk", "hjkhj"),
    @(8, '`2h',     "hjk",   "h",       "jklh",    "jk",    "hjk", "lhjk")
)

$row = 4
foreach ($student in $newStudents) {
    $students.Cells.Item($row, 1).Value = $student[0]
    $students.Cells.Item($row, 2).Value = $student[1]
    $students.Cells.Item($row, 3).Value = $student[2]
    $students.Cells.Item($row, 4).Value = $student[3]
    $students.Cells.Item($row, 5).Value = $student[4]
    $students.Cells.Item($row, 6).Value = $student[5]
    $students.Cells.Item($row, 7).Value = $student[6]
    $students.Cells.Item($row, 8).Value = $student[7]
    $row = $row + 1
}

# --- add matching password rows on "student_pswd" ---
$newPasswords = @(
    @(3, "student3", "Sam"),
    @(4, "student4", "Agnetha"),
    @(5, "student5", "Lee"),
    @(6, "student6", "Garry"),
    @(7, "student7", "Kumara"),
    @(8, "student8", '`2h')
)

$row = 4
foreach ($p in $newPasswords) {
    $pswd.Cells.Item($row, 1).Value = $p[0]
    $pswd.Cells.Item($row, 2).Value = $p[1]
    $pswd.Cells.Item($row, 3).Value = $p[2]
    $row = $row + 1
}

# update the row-count helper cell
$pswd.Range("G6").Value = 8

# the user ends up with focus on the student_pswd sheet, cell B3 selected
$pswd.Activate() | Out-Null
$pswd.Range("B3").Select() | Out-Null
